$d = $word.ActiveDocument

$replacements = @(
    @{old = "99÷9="; new = "35÷2="},
    @{old = "45÷3="; new = "75÷9="},
    @{old = "15÷7="; new = "68÷6="},
    @{old = "67÷6="; new = "59÷9="},
    @{old = "17÷8="; new = "18÷3="},
    @{old = "93÷4="; new = "84÷4="},
    @{old = "72÷3="; new = "43÷3="},
    @{old = "67÷9="; new = "81÷8="},
    @{old = "70÷4="; new = "93÷8="},
    @{old = "18÷7="; new = "31÷2="},
    @{old = "53÷8="; new = "22÷6="},
    @{old = "49÷3="; new = "27÷6="},
    @{old = "23÷7="; new = "49÷2="},
    @{old = "90÷3="; new = "92÷7="},
    @{old = "80÷7="; new = "30÷7="},
    @{old = "31÷7="; new = "68÷8="},
    @{old = "55÷8="; new = "89÷7="},
    @{old = "26÷8="; new = "63÷5="},
    @{old = "57÷3="; new = "90÷4="},
    @{old = "76÷3="; new = "22÷6="},
    @{old = "64÷5="; new = "51÷9="},
    @{old = "77÷2="; new = "72÷5="},
    @{old = "40÷6="; new = "82÷4="},
    @{old = "35÷3="; new = "98÷2="},
    @{old = "41÷9="; new = "29÷4="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
